# Update a set of imputed values in the RandomForest result data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B3"  = 6.750600000000004
    "B21" = 9.216300000000004
    "B23" = 9.00490000000001
    "B25" = 6.078799999999998
    "D27" = -8.744600000000002
    "D31" = -8.540300000000004
    "D39" = -8.1831
    "D48" = -7.4152
    "D51" = -7.737699999999998
    "D52" = -7.741399999999996
    "B53" = 5.810099999999997
    "D55" = -8.288299999999998
    "D56" = -8.0184
    "B57" = 5.262199999999998
    "D57" = -8.3642
    "B59" = 4.855699999999998
    "B69" = 5.280699999999997
    "D73" = -8.062799999999998
    "B79" = 9.411700000000003
    "B83" = 5.4985
    "D89" = -5.814400000000003
    "D90" = -8.030800000000003
    "B93" = 5.364199999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
